$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The CV-to-Job keyword match values in column F (rows 2-6) are rotated
# (shifted up by two), as part of the "top 100 Unis" / Uni CV matching update.
$ws.Range("F2").Value = "redshift : 1"
$ws.Range("F3").Value = "engineer : 1"
$ws.Range("F4").Value = "amazon : 2"
$ws.Range("F5").Value = "data engineer : 1"
$ws.Range("F6").Value = "analysis : 6"
